# Generate Report for Handback
# Updates the "Correspond Handback DateTime" column (column G) on each
# language sheet: every cell that currently holds the sheet's old
# handback-generation timestamp is refreshed to the new timestamp that
# was produced when the report was (re-)generated.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @(
    @{ Sheet = "zh-cn"; Old = "2016-02-22 08:50:29"; New = "2016-02-22 08:58:35" },
    @{ Sheet = "de-de"; Old = "2016-02-22 08:50:40"; New = "2016-02-22 08:58:45" },
    @{ Sheet = "ja-jp"; Old = "2016-02-22 08:50:50"; New = "2016-02-22 08:58:57" },
    @{ Sheet = "zh-tw"; Old = "2016-02-22 08:51:01"; New = "2016-02-22 08:59:06" }
)

foreach ($update in $sheetUpdates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 7)
        if ($cell.Text -eq $update.Old) {
            $cell.Value = $update.New
        }
    }
}
